# Applies the "Inclui projeto transporte metropolitano - '0000000'" edit:
#  - updates a handful of existing "valor_iniciativa" amounts (col D)
#  - appends a new initiative row (row 64) with a provisional instrument
#    code '0000000', description, "anexo" V, and its value
#  - tidies up the sheet view / row sizing to match what Excel leaves
#    behind after this kind of edit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Value corrections on existing rows -----------------------------------
$ws.Range("D4").Value  = 410000000     # 9288132 - Construção de pontes...
$ws.Range("D5").Value  = 887000000     # 9288133 - Recuperação de rodovias...
$ws.Range("D35").Value = 10000000      # 9288183
$ws.Range("D49").Value = 2300000       # 9288210
$ws.Range("D50").Value = 3400000       # 9288211

# --- New initiative row (row 64) -------------------------------------------
# Bring formats for B64:D64 in line with the rows just above (58-63) by
# copying their formatting before filling in the new values/format for A64.
$ws.Range("B63:D63").Copy() | Out-Null
$ws.Range("B64").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# codigo_iniciativa is a provisional instrument number "0000000" - keep it
# as text (not 0) with the same look as the other numeric columns.
$ws.Range("A64").NumberFormat = "@"
$ws.Range("A64").Value = "0000000"
$ws.Range("A64").HorizontalAlignment = -4152   # xlRight
$ws.Range("A64").VerticalAlignment = -4108     # xlCenter
$ws.Range("A64").WrapText = $true

$ws.Range("B64").Value = "Melhoria da infraestrutura dos municípios – Fortalecimento do transporte metropolitano"
$ws.Range("C64").Value = "V"
$ws.Range("D64").Value = 380000000

$ws.Rows.Item(64).RowHeight = 30.75
$ws.Rows.Item(65).RowHeight = 15

# --- View state: scroll/selection land near the newly added row -----------
$excel.ActiveWindow.ScrollRow = 60
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B64").Select() | Out-Null
